# Updated cryptos list on Wed Mar  1 21:40:50 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# every coin row, and reflects the rank swap between "Algorand" and
# "InternetComputer(DFINITY)" (rows 37/38) including their new link/price
# cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Most "Price" values look exactly like plain decimal numbers (e.g. "1.000",
# "0.9998"). Excel's normal cell-entry parsing would silently turn those into
# numeric values and drop meaningful trailing/leading zeros, so force genuine
# numeric-looking text back in as literal text (the classic leading
# apostrophe trick) while leaving already-unambiguous text (thousand-dotted
# values such as "23.547.35", volume percentages, coin names, links, …)
# untouched.
function Set-TextCell {
    param($addr, $value)

    if ($value -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $ws.Range($addr).Value = "'" + $value
    } else {
        $ws.Range($addr).Value = $value
    }
}

function Set-PriceAndVolume {
    param($row, $price, $volume)

    Set-TextCell "D$row" $price
    Set-TextCell "E$row" $volume
}

Set-PriceAndVolume 2  "23.547.35"    "  +1.53%  "
Set-PriceAndVolume 3  "1.656.30"     "  +2.81%  "
Set-PriceAndVolume 4  "1.000"        "  -0.61%  "
Set-PriceAndVolume 5  "0.9998"       "  -0.52%  "
Set-PriceAndVolume 6  "302.65"       "  +0.09%  "
Set-PriceAndVolume 7  "0.3840"       "  +1.85%  "
Set-PriceAndVolume 8  "0.3605"       "  +2.60%  "
Set-PriceAndVolume 9  "51.10"        "  -1.74%  "
Set-PriceAndVolume 10 "0.08209"      "  +2.03%  "
Set-PriceAndVolume 11 "1.238"        "  +3.60%  "
Set-PriceAndVolume 12 "1.0000"       "  -0.64%  "
Set-PriceAndVolume 13 "22.42"        "  +2.31%  "
Set-PriceAndVolume 14 "6.487"        "  +2.49%  "
Set-PriceAndVolume 15 "7.504"        "  +4.99%  "
Set-PriceAndVolume 16 "0.00001228"   "  +1.68%  "
Set-PriceAndVolume 17 "1.651.84"     "  +2.40%  "
Set-PriceAndVolume 18 "97.63"        "  +3.76%  "
Set-PriceAndVolume 19 "0.07005"      "  +1.33%  "
Set-PriceAndVolume 20 "6.793"        "  +4.85%  "
Set-PriceAndVolume 21 "17.60"        "  +2.92%  "
Set-PriceAndVolume 22 "0.9994"       "  -0.51%  "
Set-PriceAndVolume 23 "12.68"        "  +3.94%  "
Set-PriceAndVolume 24 "23.560.74"    "  +1.57%  "
Set-PriceAndVolume 25 "2.525"        "  -0.74%  "
Set-PriceAndVolume 26 "3.034"        "  -1.03%  "
Set-PriceAndVolume 27 "21.26"        "  +2.44%  "
Set-PriceAndVolume 28 "153.73"       "  +1.80%  "
Set-PriceAndVolume 29 "5.241"        "  +0.04%  "
Set-PriceAndVolume 30 "134.22"       "  +1.92%  "
Set-PriceAndVolume 31 "1.832.02"     "  +1.96%  "
Set-PriceAndVolume 32 "7.122"        "  +11.30%  "

# Row 33 (WEMIXTOKEN) only gets a refreshed volume; its price is unchanged.
Set-TextCell "E33" "  +6.12%  "

Set-PriceAndVolume 34 "12.04"        "  +6.50%  "
Set-PriceAndVolume 35 "1.062"        "  +0.15%  "
Set-PriceAndVolume 36 "0.02797"      "  +3.85%  "

# Rows 37/38 swap: "Algorand" and "InternetComputer(DFINITY)" traded ranking
# positions, each carrying its own refreshed link/price/volume.
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-PriceAndVolume 37 "6.107" "  +5.20%  "

$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-PriceAndVolume 38 "0.2500" "  +2.05%  "

Set-PriceAndVolume 39 "0.08769"      "  +1.41%  "
Set-PriceAndVolume 40 "0.07005"      "  +1.91%  "
Set-PriceAndVolume 41 "13.15"        "  +10.81%  "
Set-PriceAndVolume 42 "0.7006"       "  +2.63%  "
Set-PriceAndVolume 43 "1.338"        "  +2.52%  "
Set-PriceAndVolume 44 "15.98"        "  +5.19%  "
Set-PriceAndVolume 45 "0.6536"       "  +4.30%  "
Set-PriceAndVolume 46 "1.000"        "  -0.30%  "
Set-PriceAndVolume 47 "2.307"        "  +3.10%  "
Set-PriceAndVolume 48 "3.958"        "  +0.37%  "
Set-PriceAndVolume 49 "0.07907"      "  +0.56%  "
Set-PriceAndVolume 50 "128.21"       "  +0.41%  "
Set-PriceAndVolume 51 "1.189"        "  +2.13%  "
